$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 1929.3
$ws.Range("J17").Value2 = 269.9655
$ws.Range("L17").Value2 = 809.8965000000001
$ws.Range("N17").Value2 = -1145.8965

$ws.Range("H127").Value2 = 1333.3334

$ws.Range("H129").Value2 = 1666.0454
$ws.Range("I129").Value2 = 850.63635
$ws.Range("J129").Value2 = 2481.4546
$ws.Range("K129").Value2 = 2551.90905
$ws.Range("L129").Value2 = 7444.3638
$ws.Range("M129").Value2 = 2448.09095
$ws.Range("N129").Value2 = -17444.3638

$ws.Range("H131").Value2 = 1672.8334
$ws.Range("I131").Value2 = 1706.7273
$ws.Range("J131").Value2 = 1300
$ws.Range("K131").Value2 = 5120.1819
$ws.Range("L131").Value2 = 3900
$ws.Range("M131").Value2 = -80.18189999999959
$ws.Range("N131").Value2 = -13980

$ws.Range("H138").Value2 = 3478.12
$ws.Range("I138").Value2 = 3217.65
$ws.Range("J138").Value2 = 4520
$ws.Range("K138").Value2 = 9652.950000000001
$ws.Range("L138").Value2 = 13560
$ws.Range("M138").Value2 = -4512.950000000001
$ws.Range("N138").Value2 = -23840

$ws.Range("H141").Value2 = 4754.154
$ws.Range("I141").Value2 = 4618.091
$ws.Range("K141").Value2 = 13854.273
$ws.Range("M141").Value2 = -8674.273000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 4789.0454
$ws.Range("I32").Value2 = 5376.46
$ws.Range("J32").Value2 = 2953.375
$ws.Range("K32").Value2 = 5376.46
$ws.Range("L32").Value2 = 2953.375
$ws.Range("M32").Value2 = -5089.46
$ws.Range("N32").Value2 = -3527.375

$ws.Range("H61").Value2 = 2567.6667
$ws.Range("I61").Value2 = 1539
$ws.Range("J61").Value2 = 4625
$ws.Range("K61").Value2 = 1539
$ws.Range("L61").Value2 = 4625
$ws.Range("M61").Value2 = -1327
$ws.Range("N61").Value2 = -5049

$ws.Range("H92").Value2 = 18108.334
$ws.Range("J92").Value2 = 18108.334
$ws.Range("L92").Value2 = 18108.334
$ws.Range("N92").Value2 = -23100.334

$ws.Range("H122").Value2 = 1643.875
$ws.Range("I122").Value2 = 1441.8334
$ws.Range("J122").Value2 = 2250
$ws.Range("K122").Value2 = 4325.5002
$ws.Range("L122").Value2 = 6750
$ws.Range("M122").Value2 = -1875.5002
$ws.Range("N122").Value2 = -11650

$ws.Range("H136").Value2 = 2567.6667
$ws.Range("I136").Value2 = 1539
$ws.Range("J136").Value2 = 4625
$ws.Range("K136").Value2 = 4617
$ws.Range("L136").Value2 = 13875
$ws.Range("M136").Value2 = -2067
$ws.Range("N136").Value2 = -18975

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2331.8572
$ws.Range("I31").Value2 = 1021.9524
$ws.Range("J31").Value2 = 3314.2856
$ws.Range("K31").Value2 = 1021.9524
$ws.Range("L31").Value2 = 3314.2856
$ws.Range("M31").Value2 = -726.9524
$ws.Range("N31").Value2 = -3904.2856

$ws.Range("H34").Value2 = 2331.8572
$ws.Range("I34").Value2 = 1021.9524
$ws.Range("J34").Value2 = 3314.2856
$ws.Range("K34").Value2 = 1021.9524
$ws.Range("L34").Value2 = 3314.2856
$ws.Range("M34").Value2 = -819.9524
$ws.Range("N34").Value2 = -3718.2856

$ws.Range("H58").Value2 = 3270
$ws.Range("I58").Value2 = 2451.4
$ws.Range("J58").Value2 = 3552.276
$ws.Range("K58").Value2 = 2451.4
$ws.Range("L58").Value2 = 3552.276
$ws.Range("M58").Value2 = -2248.4
$ws.Range("N58").Value2 = -3958.276

$ws.Range("H120").Value2 = 0
$ws.Range("J120").Value2 = 0
$ws.Range("L120").Value2 = 0
$ws.Range("N120").ClearContents()

$ws.Range("H122").Value2 = 8334901.5
$ws.Range("I122").Value2 = 15626374
$ws.Range("J122").Value2 = 1790.4286
$ws.Range("K122").Value2 = 46879122
$ws.Range("L122").Value2 = 5371.2858
$ws.Range("M122").Value2 = -46876672
$ws.Range("N122").Value2 = -10271.2858

$ws.Range("H136").Value2 = 3270
$ws.Range("I136").Value2 = 2451.4
$ws.Range("J136").Value2 = 3552.276
$ws.Range("K136").Value2 = 7354.200000000001
$ws.Range("L136").Value2 = 10656.828
$ws.Range("M136").Value2 = -4804.200000000001
$ws.Range("N136").Value2 = -15756.828

$ws.Range("H141").Value2 = 8650
$ws.Range("J141").Value2 = 8650
$ws.Range("L141").Value2 = 8650
$ws.Range("N141").Value2 = -19010

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value2 = 0
$ws.Range("I87").Value2 = 0
$ws.Range("K87").Value2 = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value2 = 0
$ws.Range("I90").Value2 = 0
$ws.Range("K90").Value2 = 0
$ws.Range("M90").ClearContents()

$ws.Range("H110").Value2 = 2356
$ws.Range("I110").Value2 = 1933.3334
$ws.Range("J110").Value2 = 2990
$ws.Range("K110").Value2 = 5800.0002
$ws.Range("L110").Value2 = 8970
$ws.Range("M110").Value2 = -1710.0002
$ws.Range("N110").Value2 = -17150

$ws.Range("H114").Value2 = 1656.5
$ws.Range("I114").Value2 = 980.8333
$ws.Range("J114").Value2 = 2670
$ws.Range("K114").Value2 = 2942.4999
$ws.Range("L114").Value2 = 8010
$ws.Range("M114").Value2 = 311.5001000000002
$ws.Range("N114").Value2 = -14518

$ws.Range("H121").Value2 = 27778852
$ws.Range("I121").Value2 = 416.66666
$ws.Range("J121").Value2 = 33334540
$ws.Range("K121").Value2 = 1249.99998
$ws.Range("L121").Value2 = 100003620
$ws.Range("M121").Value2 = 60.00001999999995
$ws.Range("N121").Value2 = -100006240

$ws.Range("H131").Value2 = 1164986.9
$ws.Range("J131").Value2 = 1371578.8
$ws.Range("L131").Value2 = 4114736.4
$ws.Range("N131").Value2 = -4124816.4

$ws.Range("H133").Value2 = 6528.697
$ws.Range("J133").Value2 = 7146.7144
$ws.Range("L133").Value2 = 21440.1432
$ws.Range("N133").Value2 = -31560.1432

$ws.Range("H138").Value2 = 4747.143
$ws.Range("I138").Value2 = 5357.5
$ws.Range("J138").Value2 = 3933.3333
$ws.Range("K138").Value2 = 16072.5
$ws.Range("L138").Value2 = 11799.9999
$ws.Range("M138").Value2 = -10932.5
$ws.Range("N138").Value2 = -22079.9999

$ws.Range("H140").Value2 = 1711.9375
$ws.Range("I140").Value2 = 1623.8667
$ws.Range("J140").Value2 = 3033
$ws.Range("K140").Value2 = 4871.6001
$ws.Range("L140").Value2 = 9099
$ws.Range("M140").Value2 = 308.3999000000003
$ws.Range("N140").Value2 = -19459

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value2 = 763.0714
$ws.Range("I107").Value2 = 664.8889
$ws.Range("J107").Value2 = 939.8
$ws.Range("K107").Value2 = 664.8889
$ws.Range("L107").Value2 = 939.8
$ws.Range("M107").Value2 = 1255.1111
$ws.Range("N107").Value2 = -4779.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 2181.1177
$ws.Range("I7").Value2 = 1889.909
$ws.Range("J7").Value2 = 2715
$ws.Range("K7").Value2 = 1889.909
$ws.Range("L7").Value2 = 2715
$ws.Range("M7").Value2 = -1777.909
$ws.Range("N7").Value2 = -2939

$ws.Range("H40").Value2 = 2351
$ws.Range("I40").Value2 = 2221.6
$ws.Range("J40").Value2 = 2566.6667
$ws.Range("K40").Value2 = 2221.6
$ws.Range("L40").Value2 = 2566.6667
$ws.Range("M40").Value2 = -2085.6
$ws.Range("N40").Value2 = -2838.6667

$ws.Range("H46").Value2 = 748.25
$ws.Range("I46").Value2 = 846.2
$ws.Range("J46").Value2 = 678.2857
$ws.Range("K46").Value2 = 846.2
$ws.Range("L46").Value2 = 678.2857
$ws.Range("M46").Value2 = -658.2
$ws.Range("N46").Value2 = -1054.2857

$ws.Range("H61").Value2 = 966.3333
$ws.Range("I61").Value2 = 966.3333
$ws.Range("J61").Value2 = 0
$ws.Range("K61").Value2 = 966.3333
$ws.Range("L61").Value2 = 0
$ws.Range("M61").Value2 = -764.3333
$ws.Range("N61").ClearContents()

$ws.Range("H100").Value2 = 2000
$ws.Range("I100").Value2 = 2000
$ws.Range("K100").Value2 = 2000
$ws.Range("M100").Value2 = -1459

$ws.Range("H113").Value2 = 966.3333
$ws.Range("I113").Value2 = 966.3333
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 966.3333
$ws.Range("L113").Value2 = 0
$ws.Range("M113").Value2 = 1203.6667
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value2 = 5529.579
$ws.Range("I122").Value2 = 6689.4
$ws.Range("J122").Value2 = 3299.1538
$ws.Range("K122").Value2 = 20068.2
$ws.Range("L122").Value2 = 9897.4614
$ws.Range("M122").Value2 = -17618.2
$ws.Range("N122").Value2 = -14797.4614

$ws.Range("H126").Value2 = 2181.1177
$ws.Range("I126").Value2 = 1889.909
$ws.Range("J126").Value2 = 2715
$ws.Range("K126").Value2 = 5669.727000000001
$ws.Range("L126").Value2 = 8145
$ws.Range("M126").Value2 = -3199.727000000001
$ws.Range("N126").Value2 = -13085

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 1956.2142
$ws.Range("I122").Value2 = 1170
$ws.Range("J122").Value2 = 2742.4285
$ws.Range("K122").Value2 = 3510
$ws.Range("L122").Value2 = 8227.2855
$ws.Range("M122").Value2 = -1060
$ws.Range("N122").Value2 = -13127.2855
